$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '37.720.67'
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").Value = '2.051.15'
$ws.Range("E3").Value = '  +1.04%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("E6").Value = '  -0.73%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.41'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.377'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0838'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.21%  '

$ws.Range("E11").Value = '  -0.06%  '

$ws.Range("D12").Value = '2.354.02'
$ws.Range("E12").Value = '  +1.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.20%  '

$ws.Range("E15").Value = '  +6.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.764'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("D17").Value = '2.050.36'
$ws.Range("E17").Value = '  +0.81%  '

$ws.Range("D18").Value = '37.703.54'
$ws.Range("E18").Value = '  -0.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.90%  '

$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("E22").Value = '  -0.96%  '

$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("E24").Value = '  +1.55%  '

$ws.Range("E25").Value = '  +3.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '

$ws.Range("E28").Value = '  -0.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.82%  '

$ws.Range("E30").Value = '  -0.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.118'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("E32").Value = '  +6.72%  '

$ws.Range("E33").Value = '  -1.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.78%  '

$ws.Range("E35").Value = '  +0.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.39%  '

$ws.Range("E37").Value = '  +3.93%  '

$ws.Range("E38").Value = '  +6.18%  '

$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("E40").Value = '  +9.38%  '

$ws.Range("D41").Value = '1.525.55'
$ws.Range("E41").Value = '  +0.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.24'
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  -1.99%  '

$ws.Range("E44").Value = '  +0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0892'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.65%  '

$ws.Range("E47").Value = '  +0.08%  '

$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.26%  '

$ws.Range("D51").Value = '2.243.64'
$ws.Range("E51").Value = '  +1.12%  '
